$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add a "Deep Learning" result row at the end of the "CARBON and SILICON" section ---
# Row 24 is currently an empty (blank gap) row, so we can just fill it in directly.
$ws.Range("A24").Value = "Deep Learning"
$ws.Range("B24").Value = 87.9

# --- Add a "Deep Learning" result row at the start of the "C, Si, N, Al" section ---
# That section's data starts at row 27 (row 26 is the section header), so insert a new
# row there, pushing the existing rows (27-44) down by one.
$ws.Rows("27:27").Insert()
# Excel's Insert copies formatting from the row above (the header, style index 2);
# the new data row should be unformatted like the other plain data rows.
$ws.Cells.Item(27, 1).ClearFormats()
$ws.Range("A27").Value = "Deep Learning"
$ws.Range("B27").Value = 76.8

# --- Update the view state saved with the sheet ---
$ws.Range("C27").Select()
$win = $excel.ActiveWindow
$win.ScrollRow = 22
$win.ScrollColumn = 1
